$wb = $excel.ActiveWorkbook

$wsSrc = $wb.Worksheets.Item("DOC_SRC")
$wsExec = $wb.Worksheets.Item("Executives")

# --- Add the two new "LINK to SRC" note values with live hyperlinks ---
$wsSrc.Range("E3").Value = "http://alink.com"
$wsSrc.Range("E4").Value = "http://anotherlink.com"

$wsSrc.Hyperlinks.Add($wsSrc.Range("E3"), "http://alink.com", "", "", "http://alink.com")
$wsSrc.Hyperlinks.Add($wsSrc.Range("E4"), "http://anotherlink.com", "", "", "http://anotherlink.com")

# Keep the cells on the sheet's plain/default look (no blue underlined
# "Hyperlink" cell style) - matches the rest of the column's formatting.
$wsSrc.Range("E3").Style = "Normal"
$wsSrc.Range("E4").Style = "Normal"

# --- Row heights: every row on both sheets becomes an explicit (custom) height ---
for ($r = 1; $r -le 19; $r++) {
  $wsSrc.Rows.Item($r).RowHeight = 14.5
}

for ($r = 1; $r -le 3; $r++) {
  $wsExec.Rows.Item($r).RowHeight = 14.5
}

# --- Update the active selection on DOC_SRC ---
$wsSrc.Activate()
$wsSrc.Range("D12").Select()

Write-Output "Applied CompanyValuesAndNotes edits"
